# Apply venueCapacity (column G) values of 500 to rows 2-85,
# and update the active selection / view on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$win = $excel.ActiveWindow

# Fill column G (venueCapacity) for data rows 2 through 85 with 500
# (this also overwrites the previous G36 value of 200 with 500).
$ws.Range("G2:G85").Value = 500

# Update sheet view: scroll so row 50 is at the top, and select G71.
$win.ScrollRow = 50
$win.ScrollColumn = 1
$ws.Range("G71").Select()
